# Update test credentials on Sheet1 and leave the active selection on B2,
# matching the state captured after the edit was made in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "mngr302125"
$ws.Range("B2").Value = "pAjapEq"

$ws.Range("B2").Select()
